$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Class 1A")

# Fill column C (rows 2-20) with the name values referenced in the diff.
# Shared-string mapping introduced by the edit: "aa" (most rows), "sd" (row 5), "dss" (rows 9 and 14).
$values = @{
    2  = "aa"
    3  = "aa"
    4  = "aa"
    5  = "sd"
    6  = "aa"
    7  = "aa"
    8  = "aa"
    9  = "dss"
    10 = "aa"
    11 = "aa"
    12 = "aa"
    13 = "aa"
    14 = "dss"
    15 = "aa"
    16 = "aa"
    17 = "aa"
    18 = "aa"
    19 = "aa"
    20 = "aa"
}

foreach ($row in $values.Keys) {
    $ws.Cells.Item($row, 3).Value = $values[$row]
}

# Match the author's final active-cell selection.
$ws.Range("C14").Select()
